# Append a new data row (row 6) to the active worksheet, mirroring the
# structure/formatting of the existing rows (columns A:N).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 6

# Copy the date-formatted cell above so the new date cell picks up the
# same number format / style (rather than minting a brand-new style).
$ws.Cells.Item($row - 1, 1).Copy() | Out-Null
$ws.Cells.Item($row, 1).PasteSpecial() | Out-Null
$ws.Cells.Item($row, 1).Value = 42611.887638888889

$ws.Cells.Item($row, 2).Value = 0
$ws.Cells.Item($row, 3).Value = 47
$ws.Cells.Item($row, 4).Value = 51
$ws.Cells.Item($row, 5).Value = 100
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 10247
$ws.Cells.Item($row, 8).Value = 5879
$ws.Cells.Item($row, 9).Value = 311
$ws.Cells.Item($row, 10).Value = 51
$ws.Cells.Item($row, 11).Value = 55
$ws.Cells.Item($row, 12).Value = 1
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = "Named"
